# Add all subclasses to content list
#
# Adds 8 new subclass rows (71-78) to the "Subclasses" sheet and 1 new
# feat row (17) to the "Feats" sheet, each with a hyperlink in the
# "Source Doc" column, matching the commit "Add all subclasses to content
# list".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Subclasses sheet (sheet2) - columns:
#   A Name | B Base Class | C Revised Subclass | D Source Doc
#   E Development Status | F Release Status | G Added To Subclass Sheet?
#   H Supporting Content Status
# ---------------------------------------------------------------------
$subclasses = $wb.Worksheets.Item("Subclasses")

$newSubclassRows = @(
    @{ Row=71; A="College of Revelry";   B="Bard";      C="No";  D="Bard College - College of Revely";   E="Playtest Ready"; F="Not Released"; G="Yes"; H="None";    Url="https://editor.gmbinder.com/documents/edit/-N8RevelryCollege1" },
    @{ Row=72; A="College of Pacts";     B="Bard";      C="No";  D="Bard College - College of Pacts";    E="Playtest Ready"; F="Not Released"; G="Yes"; H="None";    Url="https://editor.gmbinder.com/documents/edit/-N8PactsCollege0001" },
    @{ Row=73; A="College of Choir";     B="Bard";      C="No";  D="Bard College - College of Choir";    E="Needs Clean Up"; F="Not Released"; G="No";  H="None";    Url="https://editor.gmbinder.com/documents/edit/-N8ChoirCollege0001"; Display="Bard College - Collge of Choir" },
    @{ Row=74; A="College of Finality";  B="Bard";      C="No";  D="Bard College - College of Finality"; E="Needs Clean Up"; F="Not Released"; G="No";  H="None";    Url="https://editor.gmbinder.com/documents/edit/-N8FinalityCollege01" },
    @{ Row=75; A="Divine Domain - War";  B="Cleric";    C="Yes"; D="Divine Domain - War";                E="Needs Clean Up"; F="Not Released"; G="No";  H="None";    Url="https://editor.gmbinder.com/documents/edit/-N8WarDomain00000001" },
    @{ Row=76; A="Demon Soul";           B="Barbarian"; C="No";  D="Path of the Demon Soul";             E="Playtest Ready"; F="Not Released"; G="Yes"; H="Unknown"; Url="https://editor.gmbinder.com/documents/edit/-N8DemonSoulPath0001"; Display="Primal Path - Demon Soul" },
    @{ Row=77; A="Warcaller";            B="Barbarian"; C="No";  D="Path of the Warcaller";              E="Playtest Ready"; F="Not Released"; G="Yes"; H="Unknown"; Url="https://editor.gmbinder.com/documents/edit/-N8WarcallerPath0001" },
    @{ Row=78; A="Storm Herald";         B="Barbarian"; C="Yes"; D="Path of the Storm Herald";           E="Playtest Ready"; F="Not Released"; G="Yes"; H="None";    Url="https://editor.gmbinder.com/documents/edit/-N8StormHeraldPath01" }
)

foreach ($r in $newSubclassRows) {
    $row = $r.Row
    $subclasses.Range("A$row").Value = $r.A
    $subclasses.Range("B$row").Value = $r.B
    $subclasses.Range("C$row").Value = $r.C
    $subclasses.Range("D$row").Value = $r.D
    $subclasses.Range("E$row").Value = $r.E
    $subclasses.Range("F$row").Value = $r.F
    $subclasses.Range("G$row").Value = $r.G
    $subclasses.Range("H$row").Value = $r.H

    if ($r.ContainsKey("Display")) {
        # Hyperlinks.Add's TextToDisplay overwrites the cell's text, so set
        # it to the (stale/typo'd) display text the original author's
        # hyperlink carries, then restore the correct cell text afterwards.
        $subclasses.Hyperlinks.Add($subclasses.Range("D$row"), $r.Url, "", "", $r.Display) | Out-Null
        $subclasses.Range("D$row").Value = $r.D
    } else {
        $subclasses.Hyperlinks.Add($subclasses.Range("D$row"), $r.Url) | Out-Null
    }
}

# Move the saved selection/scroll position to roughly where the author
# left off after adding the new rows.
$subclasses.Activate()
$subclasses.Range("E80").Select() | Out-Null

# ---------------------------------------------------------------------
# Feats sheet (sheet7) - columns:
#   A Name | B Feat Type | C Prerequisite | D Revised Feat
#   E Source Doc | F Development Status | G Release Status
# ---------------------------------------------------------------------
$feats = $wb.Worksheets.Item("Feats")

$feats.Range("A17").Value = "Warcaller"
$feats.Range("B17").Value = "Content Dependant"
$feats.Range("C17").Value = "None"
$feats.Range("D17").Value = "No"
$feats.Range("E17").Value = "Path of the Warcaller"
$feats.Range("F17").Value = "Playtest Ready"
$feats.Range("G17").Value = "Not Released"

$feats.Hyperlinks.Add($feats.Range("E17"), "https://editor.gmbinder.com/documents/edit/-N8WarcallerPath0001") | Out-Null

# Column B ("Feat Type") now needs to be wide enough to fit "Content
# Dependant" without truncating.
$feats.Columns.Item(2).AutoFit() | Out-Null

$feats.Activate()
$feats.Range("B20").Select() | Out-Null

# Leave the Subclasses tab as the active / selected sheet, matching the
# original workbook (tabSelected on Subclasses).
$subclasses.Activate()
